# Correction mineure (extraction parametre insensible a la casse)
# Splits the old "RDW-CV" data block (which erroneously contained 7 rows of
# data instead of 3) into three clean parameter blocks: RDW-CV, MicroR and
# MacroR (3 rows each: L1/L2/L3), inserting 2 net new rows into sheet "CV%".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CV%")

# Insert two blank rows right before the old row 29 so the rows that follow
# (old rows 29-119) shift down to (31-121), matching the target layout.
$ws.Rows("29:30").Insert()

# Fix up RDW-CV's L2/L3 values (were placeholders equal to L1's 0.3).
$ws.Cells.Item(27, 3).Value = 5
$ws.Cells.Item(28, 3).Value = 6

# New parameter: MicroR (L1/L2/L3)
$ws.Cells.Item(29, 1).Value = "MicroR"
$ws.Cells.Item(29, 2).Value = "L1"
$ws.Cells.Item(29, 3).Value = 0.3

$ws.Cells.Item(30, 1).Value = "MicroR"
$ws.Cells.Item(30, 2).Value = "L2"
$ws.Cells.Item(30, 3).Value = 5.03

$ws.Cells.Item(31, 1).Value = "MicroR"
$ws.Cells.Item(31, 2).Value = "L3"
$ws.Cells.Item(31, 3).Value = 8.44

# New parameter: MacroR (L1/L2/L3)
$ws.Cells.Item(32, 1).Value = "MacroR"
$ws.Cells.Item(32, 2).Value = "L1"
$ws.Cells.Item(32, 3).Value = 1.23

$ws.Cells.Item(33, 1).Value = "MacroR"
$ws.Cells.Item(33, 2).Value = "L2"
$ws.Cells.Item(33, 3).Value = 1.82

$ws.Cells.Item(34, 1).Value = "MacroR"
$ws.Cells.Item(34, 2).Value = "L3"
$ws.Cells.Item(34, 3).Value = 1.82

# Restore the view state (matches the saved workbook: scrolled a bit further
# down, selection resting on C28).
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("C28").Select()
